# Auto-generated edit script applying the Diabolos_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 52780
$ws.Range("J123").Value = 52780
$ws.Range("L123").Value = 52780
$ws.Range("N123").Value = -62580
$ws.Range("H128").Value = 54779
$ws.Range("J128").Value = 54779
$ws.Range("L128").Value = 54779
$ws.Range("N128").Value = -64739
$ws.Range("H132").Value = 2559.5557
$ws.Range("I132").Value = 2388.8076
$ws.Range("K132").Value = 7166.4228
$ws.Range("M132").Value = -4636.4228
$ws.Range("H133").Value = 293429.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 293429.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 293429.5
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -303549.5
$ws.Range("H136").Value = 60118.75
$ws.Range("J136").Value = 60118.75
$ws.Range("L136").Value = 60118.75
$ws.Range("N136").Value = -70318.75

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3141.5833
$ws.Range("I61").Value = 1494.6666
$ws.Range("J61").Value = 5886.4443
$ws.Range("K61").Value = 1494.6666
$ws.Range("L61").Value = 5886.4443
$ws.Range("M61").Value = -1282.6666
$ws.Range("N61").Value = -6310.4443
$ws.Range("H131").Value = 33332.5
$ws.Range("J131").Value = 33332.5
$ws.Range("L131").Value = 33332.5
$ws.Range("N131").Value = -43412.5
$ws.Range("H135").Value = 62854.25
$ws.Range("J135").Value = 64438.145
$ws.Range("L135").Value = 64438.145
$ws.Range("N135").Value = -74578.14499999999
$ws.Range("H136").Value = 3141.5833
$ws.Range("I136").Value = 1494.6666
$ws.Range("J136").Value = 5886.4443
$ws.Range("K136").Value = 4483.9998
$ws.Range("L136").Value = 17659.3329
$ws.Range("M136").Value = -1933.9998
$ws.Range("N136").Value = -22759.3329
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 59999
$ws.Range("J140").Value = 59999
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359
$ws.Range("H141").Value = 132499.5
$ws.Range("J141").Value = 59998.5
$ws.Range("L141").Value = 59998.5
$ws.Range("N141").Value = -70358.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 12441.167
$ws.Range("J103").Value = 12441.167
$ws.Range("L103").Value = 12441.167
$ws.Range("N103").Value = -14785.167
$ws.Range("H106").Value = 82499
$ws.Range("I106").Value = 79999
$ws.Range("J106").Value = 84999
$ws.Range("K106").Value = 79999
$ws.Range("L106").Value = 84999
$ws.Range("M106").Value = -78737
$ws.Range("N106").Value = -87523
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H132").Value = 84999
$ws.Range("J132").Value = 84999
$ws.Range("L132").Value = 84999
$ws.Range("N132").Value = -95119
$ws.Range("H137").Value = 73591.336
$ws.Range("J137").Value = 73591.336
$ws.Range("L137").Value = 73591.336
$ws.Range("N137").Value = -83791.336
$ws.Range("H138").Value = 92199.60000000001
$ws.Range("J138").Value = 92499.75
$ws.Range("L138").Value = 92499.75
$ws.Range("N138").Value = -102779.75
$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20054.25
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 26704.334
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 26704.334
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -28994.334
$ws.Range("H132").Value = 2143.1482
$ws.Range("I132").Value = 1391.7368
$ws.Range("J132").Value = 3927.75
$ws.Range("K132").Value = 4175.2104
$ws.Range("L132").Value = 11783.25
$ws.Range("M132").Value = -1645.2104
$ws.Range("N132").Value = -16843.25
$ws.Range("H134").Value = 3893.5908
$ws.Range("I134").Value = 3770
$ws.Range("K134").Value = 11310
$ws.Range("M134").Value = -8775
$ws.Range("H135").Value = 40875
$ws.Range("J135").Value = 40875
$ws.Range("L135").Value = 40875
$ws.Range("N135").Value = -51015
$ws.Range("H137").Value = 99000
$ws.Range("J137").Value = 99000
$ws.Range("L137").Value = 99000
$ws.Range("N137").Value = -109200
$ws.Range("H138").Value = 58899.8
$ws.Range("J138").Value = 58899.8
$ws.Range("L138").Value = 58899.8
$ws.Range("N138").Value = -69179.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2099.7222
$ws.Range("I34").Value = 249.66667
$ws.Range("J34").Value = 2469.7334
$ws.Range("K34").Value = 749.00001
$ws.Range("L34").Value = 7409.2002
$ws.Range("M34").Value = -665.00001
$ws.Range("N34").Value = -7577.2002
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 575
$ws.Range("I122").Value = 575
$ws.Range("K122").Value = 5175
$ws.Range("M122").Value = -2725
$ws.Range("H132").Value = 2359.9607
$ws.Range("I132").Value = 1338.3334
$ws.Range("K132").Value = 12045.0006
$ws.Range("M132").Value = -9515.000599999999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 55999
$ws.Range("J124").Value = 55999
$ws.Range("L124").Value = 55999
$ws.Range("N124").Value = -65819
$ws.Range("H130").Value = 55999
$ws.Range("J130").Value = 55999
$ws.Range("L130").Value = 55999
$ws.Range("N130").Value = -66039
$ws.Range("H133").Value = 60111
$ws.Range("J133").Value = 60111
$ws.Range("L133").Value = 60111
$ws.Range("N133").Value = -70231
$ws.Range("H137").Value = 86666.336
$ws.Range("J137").Value = 86666.336
$ws.Range("L137").Value = 86666.336
$ws.Range("N137").Value = -96866.336
$ws.Range("H138").Value = 58924.5
$ws.Range("J138").Value = 58924.5
$ws.Range("L138").Value = 58924.5
$ws.Range("N138").Value = -69204.5
$ws.Range("H140").Value = 147614.33
$ws.Range("J140").Value = 167494.5
$ws.Range("L140").Value = 167494.5
$ws.Range("N140").Value = -177854.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 40824.75
$ws.Range("I63").Value = 37249.5
$ws.Range("J63").Value = 44400
$ws.Range("K63").Value = 37249.5
$ws.Range("L63").Value = 44400
$ws.Range("M63").Value = -36500.5
$ws.Range("N63").Value = -45898
$ws.Range("H66").Value = 40824.75
$ws.Range("I66").Value = 37249.5
$ws.Range("J66").Value = 44400
$ws.Range("K66").Value = 111748.5
$ws.Range("L66").Value = 133200
$ws.Range("M66").Value = -108004.5
$ws.Range("N66").Value = -140688
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H134").Value = 18761.143
$ws.Range("J134").Value = 18761.143
$ws.Range("L134").Value = 18761.143
$ws.Range("N134").Value = -28901.143
$ws.Range("H137").Value = 59425
$ws.Range("J137").Value = 59425
$ws.Range("L137").Value = 59425
$ws.Range("N137").Value = -69625
$ws.Range("H139").Value = 47026.332
$ws.Range("I139").Value = 21650
$ws.Range("J139").Value = 59714.5
$ws.Range("K139").Value = 21650
$ws.Range("L139").Value = 59714.5
$ws.Range("M139").Value = -16510
$ws.Range("N139").Value = -69994.5
$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22227944
$ws.Range("I81").Value = 500
$ws.Range("K81").Value = 1000
$ws.Range("M81").Value = 61
$ws.Range("H84").Value = 22227944
$ws.Range("I84").Value = 500
$ws.Range("K84").Value = 5000
$ws.Range("M84").Value = 304
$ws.Range("H125").Value = 84998.5
$ws.Range("J125").Value = 84998
$ws.Range("L125").Value = 84998
$ws.Range("N125").Value = -94838
$ws.Range("H132").Value = 3666.3684
$ws.Range("I132").Value = 3518.152
$ws.Range("K132").Value = 10554.456
$ws.Range("M132").Value = -8024.456
$ws.Range("H136").Value = 5241.9375
$ws.Range("I136").Value = 5741.5654
$ws.Range("J136").Value = 3965.111
$ws.Range("K136").Value = 17224.6962
$ws.Range("L136").Value = 11895.333
$ws.Range("M136").Value = -14674.6962
$ws.Range("N136").Value = -16995.333
$ws.Range("H137").Value = 66665
$ws.Range("J137").Value = 66665
$ws.Range("L137").Value = 66665
$ws.Range("N137").Value = -76865
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

